# Commit "update file with jgit": the "Good Morning" cell (E8) is
# replaced with "GIT UPDATE", and the edited cell becomes the active
# selection, matching the sheetView <selection activeCell="E8" sqref="E8"/>
# that Excel records on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
